$d = $word.ActiveDocument

# --- Edit 1: "Block = curlybraceopen, {VarDclIdP | StatementP | (id, (VarDclWOIdAF | StatementWOIdAF | colon))}, curlybraceclose;"
# Extend VarDclWOIdAF -> VarDclWOIdAFP and StatementWOIdAF -> StatementWOIdAFP (text only, before the red "colon" run).
$d.Content.Find.Execute(", (VarDclWOIdAF | StatementWOIdAF | ", $true, $false, $false, $false, $false, $true, 1, $false, ", (VarDclWOIdAFP | StatementWOIdAFP | ", 2) | Out-Null

# Add the new "(bracketopen, (VarDclWOIdAFWOBSAF | StatementWOIdAFWOBSAF))" alternative right after the
# (red) "colon" run, before the closing "))}, ".
$d.Content.Find.Execute("))}, ", $true, $false, $false, $false, $false, $true, 1, $false, " | (bracketopen, (VarDclWOIdAFWOBSAF | StatementWOIdAFWOBSAF))))}, ", 2) | Out-Null

# --- Edit 2: "LoopStatement = (((while, ...) | (for, parenthesisopen, ((id, bracketopen, ((VariableWOIdAFWOBSAF, assignment, Expr) | VarInitDclWOIdAFWOBSAF)) | VarInitDclP | semicolon))), (Statement, Block)) | (do, ...);"
# Insert "(VarInitDclWOIdAFP | " right before the (red) "bracketopen" run, without disturbing its formatting.
# The trailing edit position of this insert is where Word's "_GoBack" bookmark ends up, so move it there too.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute("bracketopen, ((VariableWOIdAFWOBSAF") | Out-Null
$r.Collapse(1)
$insertStart = $r.Start
$insertTextA = "(VarInitDclWOIdAFP"
$insertTextB = " | "
$r.InsertBefore($insertTextA + $insertTextB)

$bmPos = $insertStart + $insertTextA.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Close the newly-opened paren and extend VarInitDclP -> VarInitDclidP.
$d.Content.Find.Execute(")) | VarInitDclP | ", $true, $false, $false, $false, $false, $true, 1, $false, "))) | VarInitDclidP | ", 2) | Out-Null
